$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values for rows 2-27 after repulling data
$values = @{
    2  = 6
    3  = 0
    4  = -6
    5  = 10
    6  = -7
    7  = 0
    8  = 3
    9  = 1
    10 = 7
    11 = 5
    12 = -2
    13 = 4
    14 = 1
    15 = -3
    16 = 8
    17 = 0
    18 = 2
    19 = 0
    20 = 1
    21 = 8
    22 = 3
    23 = 1
    24 = -2
    25 = -1
    26 = 0
    27 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
